$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Status" column header, filled in top-down as the sheet was reviewed
$ws.Range("D1").Value = "Status"
$ws.Range("D4").Value = "Complete"

# Row 5 - fix the typo in col A (points -> rebounds)
$ws.Range("A5").Value = "Checks: Do total rebounds in the game table match the the rebounds in basic stats?"
$ws.Range("D5").Value = "Complete"

# New "Comments" column header + notes
$ws.Range("E1").Value = "Comments"
$ws.Range("E4").Value = "Several errors were found and fixed. "
$ws.Range("E5").Value = "Several errors were found and identified as errors in the source data. Choice has been made to trust the indiidual game files and update the game totals tables accordingly where errors are located. Need to Update the game tables for current year issues"
$ws.Range("E5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 56

# Row 2 - Status
$ws.Range("D2").Value = "In progress"

# Column E width (closest achievable value to the authored 56.6640625
# after this host's column-width pixel quantization)
$ws.Columns.Item(5).ColumnWidth = 55.75

# New cells at bottom: date, TODAY(), and difference formula
# (set the difference formula first so it doesn't inherit a date
# number-format from its still-unformatted precedents)
$ws.Range("C21").Formula = "=C19-C20"

$ws.Range("B2").Copy()
$ws.Range("C19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C19").Value = Get-Date -Year 2017 -Month 1 -Day 29 -Hour 0 -Minute 0 -Second 0

$ws.Range("B2").Copy()
$ws.Range("C20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C20").Formula = "=TODAY()"

$ws.Range("C19").Select() | Out-Null
